$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.092903
$ws.Range("H2").Value = 0.278709
$ws.Range("I2").Value = 0.03600043090620505
$ws.Range("J2").Value = 0.03600043090620505
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.023286
$ws.Range("N2").Value = 0.069858
$ws.Range("O2").Value = 0.009310710475795457
$ws.Range("P2").Value = 0.009310710475795458
$ws.Range("Q2").Value = 0.002163339258
$ws.Range("R2").Value = 0.019470053322
$ws.Range("S2").Value = 0.0003351895891715539
$ws.Range("T2").Value = 0.0003351895891715539

$ws.Range("G3").Value = 0.092903
$ws.Range("H3").Value = 0.278709
$ws.Range("I3").Value = 0.03600043090620505
$ws.Range("J3").Value = 0.03600043090620505
$ws.Range("O3").Value = 0.05314667307834813
$ws.Range("P3").Value = 0.05314667307834814
$ws.Range("Q3").Value = 0.01234860482466666
$ws.Range("R3").Value = 0.111137443422
$ws.Range("S3").Value = 0.00191330313205174
$ws.Range("T3").Value = 0.00191330313205174

$ws.Range("G4").Value = 0.092903
$ws.Range("H4").Value = 0.278709
$ws.Range("I4").Value = 0.03600043090620505
$ws.Range("J4").Value = 0.03600043090620505
$ws.Range("M4").Value = 2.344785333333334
$ws.Range("N4").Value = 7.034356000000001
$ws.Range("O4").Value = 0.9375426164458565
$ws.Range("P4").Value = 0.9375426164458565
$ws.Range("Q4").Value = 0.2178375918226667
$ws.Range("R4").Value = 1.960538326404
$ws.Range("S4").Value = 0.03375193818498175
$ws.Range("T4").Value = 0.03375193818498175

$ws.Range("I5").Value = 0.1132051051535142
$ws.Range("J5").Value = 0.1132051051535142
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023286
$ws.Range("N5").Value = 0.069858
$ws.Range("O5").Value = 0.009310710475795457
$ws.Range("P5").Value = 0.009310710475795458
$ws.Range("Q5").Value = 0.006802725468000001
$ws.Range("R5").Value = 0.06122452921200001
$ws.Range("S5").Value = 0.001054019958466351
$ws.Range("T5").Value = 0.001054019958466351

$ws.Range("I6").Value = 0.1132051051535142
$ws.Range("J6").Value = 0.1132051051535142
$ws.Range("O6").Value = 0.05314667307834813
$ws.Range("P6").Value = 0.05314667307834814
$ws.Range("S6").Value = 0.006016474714393843
$ws.Range("T6").Value = 0.006016474714393845

$ws.Range("I7").Value = 0.1132051051535142
$ws.Range("J7").Value = 0.1132051051535142
$ws.Range("M7").Value = 2.344785333333334
$ws.Range("N7").Value = 7.034356000000001
$ws.Range("O7").Value = 0.9375426164458565
$ws.Range("P7").Value = 0.9375426164458565
$ws.Range("Q7").Value = 0.6850008977093335
$ws.Range("R7").Value = 6.165008079384001
$ws.Range("S7").Value = 0.106134610480654
$ws.Range("T7").Value = 0.106134610480654

$ws.Range("G8").Value = 2.195567
$ws.Range("H8").Value = 6.586701
$ws.Range("I8").Value = 0.8507944639402807
$ws.Range("J8").Value = 0.8507944639402808
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.023286
$ws.Range("N8").Value = 0.069858
$ws.Range("O8").Value = 0.009310710475795457
$ws.Range("P8").Value = 0.009310710475795458
$ws.Range("Q8").Value = 0.051125973162
$ws.Range("R8").Value = 0.460133758458
$ws.Range("S8").Value = 0.00792150092815755
$ws.Range("T8").Value = 0.007921500928157554

$ws.Range("G9").Value = 2.195567
$ws.Range("H9").Value = 6.586701
$ws.Range("I9").Value = 0.8507944639402807
$ws.Range("J9").Value = 0.8507944639402808
$ws.Range("O9").Value = 0.05314667307834813
$ws.Range("P9").Value = 0.05314667307834814
$ws.Range("Q9").Value = 0.2918333019286666
$ws.Range("R9").Value = 2.626499717358
$ws.Range("S9").Value = 0.04521689523190255
$ws.Range("T9").Value = 0.04521689523190256

$ws.Range("G10").Value = 2.195567
$ws.Range("H10").Value = 6.586701
$ws.Range("I10").Value = 0.8507944639402807
$ws.Range("J10").Value = 0.8507944639402808
$ws.Range("M10").Value = 2.344785333333334
$ws.Range("N10").Value = 7.034356000000001
$ws.Range("O10").Value = 0.9375426164458565
$ws.Range("P10").Value = 0.9375426164458565
$ws.Range("Q10").Value = 5.148133299950667
$ws.Range("R10").Value = 46.333199699556
$ws.Range("S10").Value = 0.7976560677802206
$ws.Range("T10").Value = 0.7976560677802207
